$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "USN"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Age"
$ws.Range("D1").Value = "Phone"

[void]$ws.Range("E1").Select()
